$wb = $excel.ActiveWorkbook

# --- Update the summary text on sheet "Hoja1" (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.43 = 13280.28 pesos`n✅ 13280.28 pesos = 3.42 = 967.4 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 291.673
$ws2.Range("O10").Value = 3873.5
$ws2.Range("N12").Value = 3884.99
$ws2.Range("O12").Value = 283
